# Auto-generated: refresh market-price derived columns (H-N) across all profession sheets.
# Source: scheduled runner data refresh (market board price snapshot update).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 499.57144
$ws.Range("I9").Value = 331.66666
$ws.Range("K9").Value = 331.66666
$ws.Range("M9").Value = -162.66666
$ws.Range("H135").Value = 1632.8334
$ws.Range("I135").Value = 759.4666999999999
$ws.Range("K135").Value = 6835.2003
$ws.Range("M135").Value = -4300.2003
$ws.Range("H138").Value = 16549.338
$ws.Range("I138").Value = 1341.3208
$ws.Range("J138").Value = 70284.336
$ws.Range("K138").Value = 4023.9624
$ws.Range("L138").Value = 210853.008
$ws.Range("M138").Value = 1116.0376
$ws.Range("N138").Value = -221133.008

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 7520.8335
$ws.Range("I61").Value = 1190.5385
$ws.Range("K61").Value = 1190.5385
$ws.Range("M61").Value = -978.5385000000001
$ws.Range("H132").Value = 1089.0555
$ws.Range("I132").Value = 859
$ws.Range("J132").Value = 5000
$ws.Range("K132").Value = 2577
$ws.Range("L132").Value = 15000
$ws.Range("M132").Value = -47
$ws.Range("N132").Value = -20060
$ws.Range("H136").Value = 7520.8335
$ws.Range("I136").Value = 1190.5385
$ws.Range("K136").Value = 3571.6155
$ws.Range("M136").Value = -1021.6155

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H128").Value = 3500
$ws.Range("I128").Value = 3500
$ws.Range("K128").Value = 10500
$ws.Range("M128").Value = -8010
$ws.Range("H134").Value = 2495.611
$ws.Range("I134").Value = 2014.3462
$ws.Range("K134").Value = 6043.0386
$ws.Range("M134").Value = -3508.0386

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4765620
$ws.Range("I31").Value = 12501209
$ws.Range("K31").Value = 12501209
$ws.Range("M31").Value = -12500914
$ws.Range("H34").Value = 4765620
$ws.Range("I34").Value = 12501209
$ws.Range("K34").Value = 12501209
$ws.Range("M34").Value = -12501007
$ws.Range("H58").Value = 19082.37
$ws.Range("I58").Value = 2087.2
$ws.Range("K58").Value = 2087.2
$ws.Range("M58").Value = -1884.2
$ws.Range("H62").Value = 5632.2
$ws.Range("I62").Value = 5995.3335
$ws.Range("K62").Value = 5995.3335
$ws.Range("M62").Value = -5371.3335
$ws.Range("H65").Value = 5632.2
$ws.Range("I65").Value = 5995.3335
$ws.Range("K65").Value = 29976.6675
$ws.Range("M65").Value = -26856.6675
$ws.Range("H86").Value = 43870.9
$ws.Range("I86").Value = 56030
$ws.Range("K86").Value = 56030
$ws.Range("M86").Value = -54907
$ws.Range("H89").Value = 43870.9
$ws.Range("I89").Value = 56030
$ws.Range("K89").Value = 280150
$ws.Range("M89").Value = -274534
$ws.Range("H99").Value = 8653.615
$ws.Range("I99").Value = 6140.7144
$ws.Range("J99").Value = 11585.333
$ws.Range("K99").Value = 6140.7144
$ws.Range("L99").Value = 11585.333
$ws.Range("M99").Value = -4642.7144
$ws.Range("N99").Value = -14581.333
$ws.Range("H126").Value = 8653.615
$ws.Range("I126").Value = 6140.7144
$ws.Range("J126").Value = 11585.333
$ws.Range("K126").Value = 18422.1432
$ws.Range("L126").Value = 34755.999
$ws.Range("M126").Value = -15952.1432
$ws.Range("N126").Value = -39695.999
$ws.Range("H132").Value = 64264.25
$ws.Range("I132").Value = 84352.414
$ws.Range("K132").Value = 253057.242
$ws.Range("M132").Value = -250527.242
$ws.Range("H134").Value = 2365.1667
$ws.Range("I134").Value = 1838.2
$ws.Range("K134").Value = 5514.6
$ws.Range("M134").Value = -2979.6
$ws.Range("H135").Value = 119997.46
$ws.Range("J135").Value = 119997.46
$ws.Range("L135").Value = 119997.46
$ws.Range("N135").Value = -130137.46
$ws.Range("H136").Value = 19082.37
$ws.Range("I136").Value = 2087.2
$ws.Range("K136").Value = 6261.599999999999
$ws.Range("M136").Value = -3711.599999999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H141").Value = 4940
$ws.Range("I141").Value = 4050
$ws.Range("J141").Value = 8500
$ws.Range("K141").Value = 12150
$ws.Range("L141").Value = 25500
$ws.Range("M141").Value = -6970
$ws.Range("N141").Value = -35860

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H15").Value = 24514.625
$ws.Range("J15").Value = 24514.625
$ws.Range("L15").Value = 24514.625
$ws.Range("N15").Value = -25090.625
$ws.Range("H74").Value = 0
$ws.Range("J74").Value = 0
$ws.Range("L74").Value = 0
$ws.Range("H77").Value = 0
$ws.Range("J77").Value = 0
$ws.Range("L77").Value = 0
$ws.Range("H80").Value = 35240.2
$ws.Range("I80").Value = 1200
$ws.Range("J80").Value = 43750.25
$ws.Range("K80").Value = 1200
$ws.Range("L80").Value = 43750.25
$ws.Range("M80").Value = -202
$ws.Range("N80").Value = -45746.25
$ws.Range("H81").Value = 24514.625
$ws.Range("J81").Value = 24514.625
$ws.Range("L81").Value = 24514.625
$ws.Range("N81").Value = -26510.625
$ws.Range("H82").Value = 0
$ws.Range("J82").Value = 0
$ws.Range("L82").Value = 0
$ws.Range("H83").Value = 35240.2
$ws.Range("I83").Value = 1200
$ws.Range("J83").Value = 43750.25
$ws.Range("K83").Value = 6000
$ws.Range("L83").Value = 218751.25
$ws.Range("M83").Value = -1008
$ws.Range("N83").Value = -228735.25
$ws.Range("H84").Value = 24514.625
$ws.Range("J84").Value = 24514.625
$ws.Range("L84").Value = 73543.875
$ws.Range("N84").Value = -83527.875
$ws.Range("H85").Value = 0
$ws.Range("J85").Value = 0
$ws.Range("L85").Value = 0
$ws.Range("H126").Value = 4166.1665
$ws.Range("I126").Value = 1999.6666
$ws.Range("K126").Value = 5998.9998
$ws.Range("M126").Value = -3528.9998
$ws.Range("H132").Value = 3963.8708
$ws.Range("I132").Value = 3880.077
$ws.Range("J132").Value = 4399.6
$ws.Range("K132").Value = 11640.231
$ws.Range("L132").Value = 13198.8
$ws.Range("M132").Value = -9110.231
$ws.Range("N132").Value = -18258.8
$ws.Range("N74").ClearContents()
$ws.Range("N77").ClearContents()
$ws.Range("N82").ClearContents()
$ws.Range("N85").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 2462.75
$ws.Range("I40").Value = 2244.8
$ws.Range("K40").Value = 2244.8
$ws.Range("M40").Value = -2108.8
$ws.Range("H122").Value = 3925.682
$ws.Range("J122").Value = 6749.2
$ws.Range("L122").Value = 20247.6
$ws.Range("N122").Value = -25147.6
$ws.Range("H132").Value = 2687.5386
$ws.Range("I132").Value = 2243.9
$ws.Range("K132").Value = 6731.700000000001
$ws.Range("M132").Value = -4201.700000000001
$ws.Range("H136").Value = 10002
$ws.Range("I136").Value = 10002
$ws.Range("K136").Value = 30006
$ws.Range("M136").Value = -27456

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 2432.2632
$ws.Range("I132").Value = 1516.6923
$ws.Range("K132").Value = 4550.0769
$ws.Range("M132").Value = -2020.0769
